# The commit removes two slides from the deck:
#   - slide 27: "CRISPR interference (CRISPRi)/CRISPR nuclease  (CRISPRn)"
#   - slide 28: "Datasets"
# After removing them, the former slides 29 ("Expected Impact") and
# 30 ("Future Work") shift up to become slides 27 and 28.

$p = $ppt.ActivePresentation

# Delete the "CRISPR interference ... / Datasets" pair of slides.
# Deleting slide 27 twice removes both, because after the first delete
# the old slide 28 ("Datasets") becomes the new slide 27.
$p.Slides.Item(27).Delete()
$p.Slides.Item(27).Delete()
